# The workbook's weekly price list gets a new week's record inserted at the
# top of the data block (row 31), pushing all the later rows (old 31-48)
# down by one (to 32-49). The new record reuses the same fixed/static
# attributes (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria,
# Variedad, Calidad, Unidad de comercializacion, Origen, Kg o Unidades,
# Clasificacion) as the rest of the "Jengibre" table, only the
# date/volume/price columns (D, J, K, L, M, P) are new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 31; this shifts existing rows 31..48
# down to 32..49 and grows the sheet dimension to A1:R49.
$ws.Rows("31:31").Insert()

# The cells that just shifted down into row 32 still hold the values that,
# before the edit, lived in row 31 - reuse them for the columns that are
# identical across every record in this table.
$ws.Range("A31").Value = $ws.Range("A32").Value()
$ws.Range("B31").Value = $ws.Range("B32").Value()
$ws.Range("C31").Value = $ws.Range("C32").Value()
$ws.Range("E31").Value = $ws.Range("E32").Value()
$ws.Range("F31").Value = $ws.Range("F32").Value()
$ws.Range("G31").Value = $ws.Range("G32").Value()
$ws.Range("H31").Value = $ws.Range("H32").Value()
$ws.Range("I31").Value = $ws.Range("I32").Value()
$ws.Range("N31").Value = $ws.Range("N32").Value()
$ws.Range("O31").Value = $ws.Range("O32").Value()
$ws.Range("Q31").Value = $ws.Range("Q32").Value()
$ws.Range("R31").Value = $ws.Range("R32").Value()

# New week's specific data for the inserted row.
$ws.Range("D31").Value = 44466
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 12000
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = 12400
$ws.Range("P31").Value = 954
